$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.NumberFormat = "General"
}

Set-TextValue "D2" "327.08"
Set-TextValue "E2" "-1.25%"
Set-TextValue "D3" "43.71"
Set-TextValue "D4" "5.474"
Set-TextValue "E4" "-3.63%"
Set-TextValue "D5" "0.08075"
Set-TextValue "E5" "-4.00%"
Set-TextValue "D6" "8.656"
Set-TextValue "E6" "-1.86%"
Set-TextValue "E7" "-4.31%"
Set-TextValue "D8" "1.881"
Set-TextValue "E8" "-5.51%"
Set-TextValue "D9" "2.700"
Set-TextValue "E9" "-8.10%"
Set-TextValue "D10" "0.9372"
Set-TextValue "E10" "1.19%"
Set-TextValue "D11" "0.1192"
Set-TextValue "E11" "-6.86%"
Set-TextValue "D12" "0.1908"
Set-TextValue "E12" "-2.84%"
Set-TextValue "D13" "0.09574"
Set-TextValue "E13" "2.39%"
Set-TextValue "D14" "0.04090"
Set-TextValue "E14" "2.24%"
Set-TextValue "E15" "0.53%"
Set-TextValue "D16" "0.001272"
Set-TextValue "E16" "-2.51%"
Set-TextValue "D17" "0.005995"
Set-TextValue "E17" "-2.06%"
Set-TextValue "D18" "3.576"
Set-TextValue "E18" "4.40%"
Set-TextValue "D20" "8.632"
Set-TextValue "E20" "-3.82%"
Set-TextValue "D21" "0.1367"
Set-TextValue "E21" "0.07%"
Set-TextValue "E22" "-0.77%"
Set-TextValue "D23" "0.04352"
Set-TextValue "E23" "-1.49%"
Set-TextValue "D24" "0.001235"
Set-TextValue "E24" "-0.87%"
Set-TextValue "D25" "0.004323"
Set-TextValue "E25" "-0.94%"
Set-TextValue "D26" "0.0001234"
Set-TextValue "E26" "3.59%"
Set-TextValue "D27" "0.0004001"
Set-TextValue "E27" "0.03%"
Set-TextValue "D39" "0.02662"
Set-TextValue "E39" "-6.00%"
Set-TextValue "D40" "0.05420"
Set-TextValue "E40" "-1.87%"
Set-TextValue "D41" "0.007654"
Set-TextValue "E41" "-3.15%"
Set-TextValue "D42" "0.01015"
Set-TextValue "E42" "13.04%"
Set-TextValue "D43" "0.1388"
Set-TextValue "E43" "-3.49%"
Set-TextValue "D44" "0.002092"
Set-TextValue "E44" "0.47%"
Set-TextValue "D45" "0.009889"
Set-TextValue "E45" "-10.88%"
Set-TextValue "D46" "0.00006877"
Set-TextValue "E46" "-1.06%"
Set-TextValue "D47" "0.00000000752"
Set-TextValue "E47" "0.03%"
Set-TextValue "D48" "0.003481"
Set-TextValue "E48" "5.67%"
Set-TextValue "D49" "0.002277"
Set-TextValue "E49" "-0.30%"
Set-TextValue "D50" "0.00002107"
Set-TextValue "E50" "0.03%"
Set-TextValue "D51" "0.0002006"
Set-TextValue "E51" "0.03%"
